$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Round 1 "Roots and Relics" bracket winners (column O) for rows 19-33
# (order chosen so new shared-string entries are appended in the same
# sequence as the canonical edit: GinkoRhizo, Ceolacanth, Saimaa Ringed Seal)
$ws.Range("O19").Value = "GinkoRhizo"
$ws.Range("O25").Value = "Ceolacanth"
$ws.Range("O23").Value = "Saimaa Ringed Seal"
$ws.Range("O21").Value = "Asian Forest Tortoise"
$ws.Range("O27").Value = "Horseshoe Crab"
$ws.Range("O29").Value = "Alligator Gar"
$ws.Range("O31").Value = "Tuatara"
$ws.Range("O33").Value = "Frilled Shark"

# Update the active selection to match the author's final cursor position
$ws.Range("O21").Select()
